# Natmi following Dr Hou advice
# Update LR-pair edge statistics for the Ifnb1-Ifnar1 interaction now that
# the ligand/receptor expressing-cell counts changed from 1 to 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  E = 3; G = 0.290242; H = 0.870726; K = 3;
       M = 25.17918166666666;  N = 75.53754499999999;
       O = 0.2628209717145306; P = 0.2628209717145306;
       Q = 7.308056045296667;  R = 65.77250440767;
       S = 0.2628209717145306; T = 0.2628209717145306 },
    @{ Row = 3;  E = 3; G = 0.290242; H = 0.870726; K = 3;
       M = 20.25845733333333;  N = 60.775372;
       O = 0.2114583195065722; P = 0.2114583195065722;
       Q = 5.879855173341333;  R = 52.918696560072;
       S = 0.2114583195065722; T = 0.2114583195065722 },
    @{ Row = 4;  E = 3; G = 0.290242; H = 0.870726; K = 3;
       M = 38.16548666666667;  N = 114.49646;
       O = 0.3983723706545385; P = 0.3983723706545386;
       Q = 11.07722718110667;  R = 99.69504462995999;
       S = 0.3983723706545385; T = 0.3983723706545386 },
    @{ Row = 5;  E = 3; G = 0.290242; H = 0.870726; K = 3;
       M = 12.20042266666667;  N = 36.601268;
       O = 0.1273483381243586; P = 0.1273483381243586;
       Q = 3.541075075618667;  R = 31.869675680568;
       S = 0.1273483381243586; T = 0.1273483381243586 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("E$row").Value = $r.E
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("K$row").Value = $r.K
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
    $ws.Range("R$row").Value = $r.R
    $ws.Range("S$row").Value = $r.S
    $ws.Range("T$row").Value = $r.T
}
